$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1. Swap B13 / B14 text content ("Conteudo das linhas" <-> "Conteudo das linhas (em ordem)")
$tmp = $ws.Range('B13').Value2
$ws.Range('B13').Value2 = $ws.Range('B14').Value2
$ws.Range('B14').Value2 = $tmp

# 2. Remove the two pictures/shapes (logos) from the worksheet
$shapeCount = $ws.Shapes.Count
for ($i = $shapeCount; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# 3. Update the view: selection moves to C12, top-left cell resets to default (A1)
$ws.Range('C12').Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
